$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new D value (price text), new E value (volume % text).
# $null for D means that column is unchanged for that row.
$updates = @(
    @{ Row = 2;  D = "28.503.34"; E = "  +0.38%  " },
    @{ Row = 3;  D = "1.569.73";  E = "  -1.69%  " },
    @{ Row = 4;  D = $null;       E = "  +0.07%  " },
    @{ Row = 5;  D = "212.12";    E = "  -1.23%  " },
    @{ Row = 6;  D = $null;       E = "  -0.78%  " },
    @{ Row = 7;  D = $null;       E = "  +0.09%  " },
    @{ Row = 8;  D = "46.15";     E = "  +5.04%  " },
    @{ Row = 9;  D = "24.02";     E = "  -0.29%  " },
    @{ Row = 10; D = $null;       E = "  -1.79%  " },
    @{ Row = 11; D = $null;       E = "  -1.90%  " },
    @{ Row = 12; D = "0.0888";    E = "  -0.01%  " },
    @{ Row = 13; D = "1.795.75";  E = "  -1.57%  " },
    @{ Row = 14; D = "1.570.84";  E = "  -1.58%  " },
    @{ Row = 15; D = "0.521";     E = "  -2.27%  " },
    @{ Row = 16; D = "28.482.00"; E = "  +0.27%  " },
    @{ Row = 17; D = $null;       E = "  -2.54%  " },
    @{ Row = 18; D = $null;       E = "  -1.77%  " },
    @{ Row = 19; D = "230.56";    E = "  +0.83%  " },
    @{ Row = 20; D = $null;       E = "  -2.24%  " },
    @{ Row = 21; D = "0.0₃0691";  E = "  -2.93%  " },
    @{ Row = 22; D = $null;       E = "  +0.07%  " },
    @{ Row = 23; D = $null;       E = "  -5.92%  " },
    @{ Row = 24; D = "9.09";      E = "  -2.58%  " },
    @{ Row = 25; D = $null;       E = "  +7.61%  " },
    @{ Row = 26; D = "150.92";    E = "  -0.71%  " },
    @{ Row = 27; D = $null;       E = "  -1.41%  " },
    @{ Row = 28; D = $null;       E = "  -2.83%  " },
    @{ Row = 29; D = $null;       E = "  -3.86%  " },
    @{ Row = 30; D = $null;       E = "  +0.07%  " },
    @{ Row = 31; D = "0.0481";    E = "  +0.93%  " },
    @{ Row = 32; D = $null;       E = "  -3.58%  " },
    @{ Row = 33; D = $null;       E = "  -1.52%  " },
    @{ Row = 34; D = $null;       E = "  -2.19%  " },
    @{ Row = 35; D = "1.393.98";  E = "  -0.40%  " },
    @{ Row = 36; D = $null;       E = "  +0.69%  " },
    @{ Row = 37; D = $null;       E = "  -4.04%  " },
    @{ Row = 38; D = $null;       E = "  +0.54%  " },
    @{ Row = 39; D = $null;       E = "  +3.69%  " },
    @{ Row = 40; D = $null;       E = "  -1.48%  " },
    @{ Row = 41; D = $null;       E = "  -3.85%  " },
    @{ Row = 42; D = $null;       E = "  +0.09%  " },
    @{ Row = 43; D = $null;       E = "  -0.56%  " },
    @{ Row = 44; D = $null;       E = "  -3.73%  " },
    @{ Row = 45; D = $null;       E = "  +2.26%  " },
    @{ Row = 46; D = $null;       E = "  -5.03%  " },
    @{ Row = 47; D = $null;       E = "  -1.34%  " },
    @{ Row = 48; D = "62.83";     E = "  -2.65%  " },
    @{ Row = 49; D = "1.706.74";  E = "  -1.47%  " },
    @{ Row = 50; D = "86.20";     E = "  -1.70%  " },
    @{ Row = 51; D = $null;       E = "  -1.68%  " }
)

# Rows whose new D value reads as a plain number (e.g. "212.12"). Excel's
# Range.Value setter auto-converts numeric-looking text to a Number, but
# the source column stores these as text (e.g. "212.05", "0.520" keep
# trailing zeros). Pre-formatting those specific cells as Text before the
# assignment keeps them as real text values, same as the rest of the column;
# resetting the style to "Normal" afterwards drops the temporary Text
# number-format again so the cell's style matches the untouched cells.
$textCoerceRows = @(5, 8, 9, 12, 15, 19, 24, 26, 31, 48, 50)
foreach ($r in $textCoerceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

foreach ($r in $textCoerceRows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}
